$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency price / volume(1h) snapshot.
# Row 28/29 also swap (Monero <-> Filecoin changed rank order).

# Column D ("Price") cells hold plain numeric-looking text (e.g. "1.648.17")
# in the source workbook, so force Text format before writing so Excel
# does not reinterpret them as numbers.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D16","D17","D18","D19","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = [ordered]@{
    'D2' = '24.201.69'
    'E2' = '  -2.99%  '
    'D3' = '1.648.17'
    'E3' = '  -3.16%  '
    'D4' = '1.003'
    'E4' = '  +0.00%  '
    'D5' = '308.59'
    'E5' = '  -2.32%  '
    'D6' = '1.002'
    'E6' = '  +0.01%  '
    'D7' = '0.3904'
    'E7' = '  -1.70%  '
    'D8' = '0.3877'
    'E8' = '  -3.73%  '
    'D9' = '1.002'
    'E9' = '  -0.12%  '
    'D10' = '1.369'
    'E10' = '  -6.82%  '
    'D11' = '49.06'
    'E11' = '  -7.00%  '
    'D12' = '0.08488'
    'E12' = '  -3.62%  '
    'D13' = '24.34'
    'E13' = '  -6.19%  '
    'D14' = '7.190'
    'E14' = '  -3.61%  '
    'E15' = '  -4.20%  '
    'D16' = '7.540'
    'E16' = '  -5.37%  '
    'D17' = '1.647.35'
    'E17' = '  -3.89%  '
    'D18' = '95.01'
    'E18' = '  -1.31%  '
    'D19' = '0.06947'
    'E19' = '  -3.41%  '
    'D20' = '21.19'
    'E20' = '  +2.95%  '
    'D21' = '6.979'
    'E21' = '  -5.01%  '
    'E22' = '  -0.05%  '
    'D23' = '13.87'
    'E23' = '  -3.93%  '
    'D24' = '24.212.06'
    'E24' = '  -2.99%  '
    'D25' = '2.338'
    'E25' = '  -0.55%  '
    'D26' = '2.757'
    'E26' = '  -6.55%  '
    'D27' = '22.62'
    'E27' = '  -4.70%  '
    'B28' = 'Filecoin'
    'C28' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D28' = '8.802'
    'E28' = '  +5.59%  '
    'B29' = 'Monero'
    'C29' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D29' = '158.21'
    'E29' = '  -1.99%  '
    'D30' = '143.05'
    'E30' = '  -4.40%  '
    'D31' = '5.402'
    'E31' = '  -12.81%  '
    'D32' = '2.450'
    'E32' = '  -6.79%  '
    'D33' = '1.828.51'
    'E33' = '  -3.84%  '
    'D34' = '7.081'
    'E34' = '  -1.70%  '
    'D35' = '0.08111'
    'E35' = '  -5.28%  '
    'D36' = '0.9932'
    'E36' = '  -5.27%  '
    'D37' = '0.02957'
    'E37' = '  -6.07%  '
    'D38' = '0.2714'
    'E38' = '  -5.19%  '
    'D39' = '0.09323'
    'E39' = '  -2.33%  '
    'D40' = '1.482'
    'E40' = '  -0.02%  '
    'D41' = '10.07'
    'E41' = '  -7.57%  '
    'D42' = '0.7668'
    'E42' = '  -7.11%  '
    'D43' = '13.20'
    'E43' = '  -5.57%  '
    'D44' = '16.15'
    'E44' = '  -6.64%  '
    'D45' = '2.505'
    'E45' = '  -6.74%  '
    'D46' = '0.6917'
    'E46' = '  -6.40%  '
    'D47' = '4.097'
    'E47' = '  -3.75%  '
    'E48' = '  -0.11%  '
    'D49' = '0.08459'
    'E49' = '  -3.29%  '
    'D50' = '1.270'
    'E50' = '  -9.72%  '
    'D51' = '134.60'
    'E51' = '  -3.20%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
